$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NEC Table 9 data")

$ws.Range("A12").Value = "1/0"
$ws.Range("A13").Value = "2/0"
$ws.Range("A14").Value = "3/0"
$ws.Range("A15").Value = "4/0"

$ws.Range("A31").Select()
